$d = $word.ActiveDocument
$d.Content.Find.Execute("692÷4=173, 0", $true, $false, $false, $false, $false, $true, 1, $false, "753÷4=188, 1", 2) | Out-Null
$d.Content.Find.Execute("335÷7=47, 6", $true, $false, $false, $false, $false, $true, 1, $false, "804÷5=160, 4", 2) | Out-Null
$d.Content.Find.Execute("844÷7=120, 4", $true, $false, $false, $false, $false, $true, 1, $false, "541÷5=108, 1", 2) | Out-Null
$d.Content.Find.Execute("804÷9=89, 3", $true, $false, $false, $false, $false, $true, 1, $false, "202÷4=50, 2", 2) | Out-Null
$d.Content.Find.Execute("835÷6=139, 1", $true, $false, $false, $false, $false, $true, 1, $false, "811÷6=135, 1", 2) | Out-Null
$d.Content.Find.Execute("838÷8=104, 6", $true, $false, $false, $false, $false, $true, 1, $false, "870÷4=217, 2", 2) | Out-Null
$d.Content.Find.Execute("539÷4=134, 3", $true, $false, $false, $false, $false, $true, 1, $false, "802÷4=200, 2", 2) | Out-Null
$d.Content.Find.Execute("866÷5=173, 1", $true, $false, $false, $false, $false, $true, 1, $false, "434÷2=217, 0", 2) | Out-Null
$d.Content.Find.Execute("521÷2=260, 1", $true, $false, $false, $false, $false, $true, 1, $false, "967÷9=107, 4", 2) | Out-Null
$d.Content.Find.Execute("204÷2=102, 0", $true, $false, $false, $false, $false, $true, 1, $false, "197÷6=32, 5", 2) | Out-Null
$d.Content.Find.Execute("555÷7=79, 2", $true, $false, $false, $false, $false, $true, 1, $false, "575÷8=71, 7", 2) | Out-Null
$d.Content.Find.Execute("276÷8=34, 4", $true, $false, $false, $false, $false, $true, 1, $false, "469÷4=117, 1", 2) | Out-Null
$d.Content.Find.Execute("328÷3=109, 1", $true, $false, $false, $false, $false, $true, 1, $false, "564÷2=282, 0", 2) | Out-Null
$d.Content.Find.Execute("625÷5=125, 0", $true, $false, $false, $false, $false, $true, 1, $false, "833÷9=92, 5", 2) | Out-Null
$d.Content.Find.Execute("562÷7=80, 2", $true, $false, $false, $false, $false, $true, 1, $false, "848÷6=141, 2", 2) | Out-Null
$d.Content.Find.Execute("948÷9=105, 3", $true, $false, $false, $false, $false, $true, 1, $false, "216÷3=72, 0", 2) | Out-Null
$d.Content.Find.Execute("267÷9=29, 6", $true, $false, $false, $false, $false, $true, 1, $false, "234÷8=29, 2", 2) | Out-Null
$d.Content.Find.Execute("710÷8=88, 6", $true, $false, $false, $false, $false, $true, 1, $false, "907÷7=129, 4", 2) | Out-Null
$d.Content.Find.Execute("734÷7=104, 6", $true, $false, $false, $false, $false, $true, 1, $false, "254÷2=127, 0", 2) | Out-Null
$d.Content.Find.Execute("265÷6=44, 1", $true, $false, $false, $false, $false, $true, 1, $false, "438÷4=109, 2", 2) | Out-Null
$d.Content.Find.Execute("563÷9=62, 5", $true, $false, $false, $false, $false, $true, 1, $false, "749÷3=249, 2", 2) | Out-Null
$d.Content.Find.Execute("769÷8=96, 1", $true, $false, $false, $false, $false, $true, 1, $false, "668÷4=167, 0", 2) | Out-Null
$d.Content.Find.Execute("276÷7=39, 3", $true, $false, $false, $false, $false, $true, 1, $false, "586÷5=117, 1", 2) | Out-Null
$d.Content.Find.Execute("963÷3=321, 0", $true, $false, $false, $false, $false, $true, 1, $false, "235÷5=47, 0", 2) | Out-Null
$d.Content.Find.Execute("641÷4=160, 1", $true, $false, $false, $false, $false, $true, 1, $false, "629÷9=69, 8", 2) | Out-Null
